$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reference an unstyled data cell so we can restore default style after forcing text format
$defaultStyle = $ws.Range("B2").Style

$ws.Range("D2").Value = '64.507.35'
$ws.Range("E2").Value = '  +5.40%  '

$ws.Range("D3").Value = '3.079.21'
$ws.Range("E3").Value = '  +3.43%  '

$ws.Range("E4").Value = '  +0.04%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '554.71'
$ws.Range("D5").Style = $defaultStyle
$ws.Range("E5").Value = '  +2.08%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '143.40'
$ws.Range("D6").Style = $defaultStyle
$ws.Range("E6").Value = '  +10.24%  '

$ws.Range("E7").Value = '  -0.03%  '

$ws.Range("D8").Value = '3.075.29'
$ws.Range("E8").Value = '  +3.40%  '

$ws.Range("E9").Value = '  +0.79%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.91'
$ws.Range("D10").Style = $defaultStyle
$ws.Range("E10").Value = '  +16.73%  '

$ws.Range("E11").Value = '  +6.18%  '

$ws.Range("E12").Value = '  +3.40%  '

$ws.Range("B13").Value = 'Avalanche'
$ws.Range("C13").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '35.20'
$ws.Range("D13").Style = $defaultStyle
$ws.Range("E13").Value = '  +5.48%  '

$ws.Range("B14").Value = 'ShibaInu'
$ws.Range("C14").Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000226'
$ws.Range("D14").Style = $defaultStyle
$ws.Range("E14").Value = '  +4.22%  '

$ws.Range("D15").Value = '3.577.87'

$ws.Range("D16").Value = '64.476.11'
$ws.Range("E16").Value = '  +5.40%  '

$ws.Range("D17").Value = '3.082.45'
$ws.Range("E17").Value = '  +3.63%  '

$ws.Range("E18").Value = '  -0.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.71'
$ws.Range("D19").Style = $defaultStyle
$ws.Range("E19").Value = '  +2.89%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '480.09'
$ws.Range("D20").Style = $defaultStyle
$ws.Range("E20").Value = '  +2.07%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '13.61'
$ws.Range("D21").Style = $defaultStyle
$ws.Range("E21").Value = '  +4.99%  '

$ws.Range("E22").Value = '  +2.05%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.52'
$ws.Range("D23").Style = $defaultStyle
$ws.Range("E23").Value = '  +9.11%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '13.11'
$ws.Range("D24").Style = $defaultStyle
$ws.Range("E24").Value = '  +10.57%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '80.74'
$ws.Range("D25").Style = $defaultStyle
$ws.Range("E25").Value = '  +1.87%  '

$ws.Range("E26").Value = '  +0.54%  '

$ws.Range("E27").Value = '  +4.65%  '

$ws.Range("E28").Value = '  +4.77%  '

$ws.Range("E29").Value = '  +9.92%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.999'
$ws.Range("D30").Style = $defaultStyle
$ws.Range("E30").Value = '  +0.04%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '26.04'
$ws.Range("D31").Style = $defaultStyle
$ws.Range("E31").Value = '  +3.06%  '

$ws.Range("E32").Value = '  +3.00%  '

$ws.Range("E33").Value = '  +7.05%  '

$ws.Range("E34").Value = '  +6.00%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '55.22'
$ws.Range("D35").Style = $defaultStyle
$ws.Range("E35").Value = '  +1.50%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.07'
$ws.Range("D36").Style = $defaultStyle
$ws.Range("E36").Value = '  +4.58%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '464.29'
$ws.Range("D37").Style = $defaultStyle
$ws.Range("E37").Value = '  +4.86%  '

$ws.Range("E38").Value = '  +9.21%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0824'
$ws.Range("D39").Style = $defaultStyle
$ws.Range("E39").Value = '  +5.48%  '

$ws.Range("D40").Value = '3.021.62'
$ws.Range("E40").Value = '  -2.85%  '

$ws.Range("E41").Value = '  +2.18%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '8.26'
$ws.Range("D42").Style = $defaultStyle
$ws.Range("E42").Value = '  +3.50%  '

$ws.Range("E43").Value = '  +18.27%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '27.79'
$ws.Range("D44").Style = $defaultStyle
$ws.Range("E44").Value = '  +11.00%  '

$ws.Range("E45").Value = '  +8.16%  '

$ws.Range("E47").Value = '  +7.89%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.112'
$ws.Range("D48").Style = $defaultStyle

$ws.Range("B49").Value = 'Monero'
$ws.Range("C49").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '117.44'
$ws.Range("D49").Style = $defaultStyle
$ws.Range("E49").Value = '  +3.57%  '

$ws.Range("B50").Value = 'PEPE'
$ws.Range("C50").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D50").Value = '0.0₃0513'
$ws.Range("E50").Value = '  +8.46%  '

$ws.Range("E51").Value = '  +5.07%  '
